$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Objetivos:) - value was "Aulas expositivas..." -> now a professor entry
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13 (Programa resumido:) - value was "60 h" -> now an activation date.
# Copy an existing text cell that already holds "01/01/2023" as a shared
# string and paste its value only, so Excel doesn't reinterpret the text
# as a date serial number.
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row 15 (Programa:) - value was "Semestral" -> now a professor entry
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

# Row 18 (Método:) - value was "01/01/2023" -> now a different professor entry
$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"
